$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New meeting rows (row 3: Dec 20 2024, row 4: Dec 25 2024) ---
# Row 3
$ws.Range("B3").Value = 45646
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "Reviewed Queries, Models, and Project Load Ahead"
$ws.Range("D3").Value = "ma, is, se, cl, ce"

# Row 4
$ws.Range("B4").Value = 45651
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Value = "Reviewed Model & UI, Planned to write more on Overleaf & Get An MVP "

$excel.CutCopyMode = $false

# --- Widen the Topics column (C) to fit the new, longer text ---
$ws.Columns.Item(3).ColumnWidth = 62.7265625

# --- Move the active selection ---
$ws.Range("C10").Select() | Out-Null

# --- Window geometry (best effort; matches the author's resized window) ---
$win = $wb.Windows.Item(1)
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 11020
